$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56)

$bvals = @(
    'P8228'
    'P3709'
    'P8220'
    'P4259'
    'P0TMP092'
    'P8182'
    'P8210'
    'P8268'
    'P4242'
    'P8217'
    'P8263'
    'P0TMP080'
    'P8205'
    'P0TMPT007'
    'P8213'
    'P8273'
    'P5651'
    'P3214'
    'P8267'
    '?'
    'https://lod.dila.edu.tw/resource.php?id=A000089'
    'P2548'
    'P4263'
    'P2956'
    'P753'
    'P8261'
    'P8266'
    'P3379'
    'P2637'
    'P4CZ15137'
    'P8151'
    'P8249'
    'P0TMP104'
    'P8211'
    'P8206'
    'P8171'
    'P8245'
    'P00KG07267'
    'P1KG8854'
    'P3456'
    'P8093'
    'P8269'
    'P8209'
    'P4258'
    'P0RK8'
    'P4CZ16780'
    'P8183'
    'P4CZ16819'
    'P4255'
    'P3285'
    'P0TMP098'
    'P8260'
    'P8222'
    'P8265'
    'P8219'
)

$cvals = @(
    '{''eft:surendrabodhi''}'
    '{''eft:phakpa-sherab''}'
    '{''eft:devacandra''}'
    '{''eft:dpal-gyi-lhun-po'', ''eft:palgyi-lh-npo'', ''eft:ban-de-dpal-gyi-lhun-po''}'
    '{''eft:anandasri-s-''}'
    '{''eft:paltsek'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:dpal-brtsegs'', ''eft:ban-de-dpal-brtsegs'', ''eft:ska-ba-dpal-brtsegs''}'
    '{''eft:danasila''}'
    '{''eft:buddhaprabha''}'
    '{''eft:sherab-lekpa''}'
    '{''eft:jnanagarbha'', ''eft:t-jnanagarbha''}'
    '{''eft:leki-d-''}'
    '{''eft:hwa-shang-zab-mo''}'
    '{''eft:band-yesh-de'', ''eft:yesh-d-'', ''eft:band-yesh-d-'', ''eft:zhang-yesh-d-'', ''eft:yesh-d-ye-shes-sde-'', ''eft:ye-shes-sde''}'
    '{''eft:rnam-par-mi-rtog-pa''}'
    '{''eft:t-vidyakarasimha'', ''eft:vidyakarasimha''}'
    '{''eft:rinchen-tso'', ''eft:rin-chen-tsho''}'
    '{''eft:pa-tshab-nyi-ma-grags''}'
    '{''eft:danasila''}'
    '{''eft:vijayasila''}'
    '{''eft:sakyasena''}'
    '{''eft:siladharma''}'
    '{''eft:prajnavarma'', ''eft:prajnavarman''}'
    '{''eft:dge-ba-dpal''}'
    '{''eft:krsnapandita''}'
    '{''eft:rin-chen-bzang-po''}'
    '{''eft:munivarman'', ''eft:munivarma''}'
    '{''eft:dharmatasila'', ''eft:ch-nyi-tsultrim''}'
    '{''eft:dipamkarasrijnana'', ''eft:dipamkara-srijnana''}'
    '{''eft:trakpa-gyaltsen''}'
    '{''eft:kumarakalasa''}'
    '{''eft:gayadhara''}'
    '{''eft:dharmakara''}'
    '{''eft:punyasambhava''}'
    '{''eft:vidyakaraprabha''}'
    '{''eft:celu''}'
    '{''eft:dharmasribhadra''}'
    '{''eft:buddhakaravarma''}'
    '{''eft:sarvajnadeva'', ''eft:sarvanyadeva''}'
    '{''eft:silendrabodhi'', ''eft:srilendrabodhi'', ''eft:surendrabodhi''}'
    '{''eft:tshul-khrims-rgyal-ba''}'
    '{''eft:kamalagupta''}'
    '{''eft:dgon-gling-rma''}'
    '{''eft:dzi-na-mi-tra-k-'', ''eft:jinamitra'', ''eft:jinamitra-k-''}'
    '{''eft:dpal-byor''}'
    '{''eft:dharmapala''}'
    '{''eft:manjusrigarbha''}'
    '{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}'
    '{''eft:sakyaprabha''}'
    '{''eft:ye-shes-snying-po'', ''eft:t-jnanagarbha'', ''eft:yesh-nyingpo''}'
    '{''eft:sakya-yesh-''}'
    '{''eft:jinavara''}'
    '{''eft:dpal-dbyangs''}'
    '{''eft:jnanasiddhi'', ''eft:jnanasidhi''}'
    '{''eft:ratnaraksita''}'
    '{''eft:visuddhasimha''}'
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
}
